$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92-167 down to 93-168.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new data record.
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 44586
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100112039
$ws.Cells.Item(92, 7).Value = "Ciboulette"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 240
$ws.Cells.Item(92, 11).Value = 3000
$ws.Cells.Item(92, 12).Value = 3500
$ws.Cells.Item(92, 13).Value = 3250
$ws.Cells.Item(92, 14).Value = "`$/docena de atados"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 1083
$ws.Cells.Item(92, 17).Value = 3
$ws.Cells.Item(92, 18).Value = "Hortaliza"
